$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" / "_new" header suffixes to "_FV2310" / "_FV2404"
#    A1:J1 hold the "_old" headers, L1:U1 hold the "_new" headers (K1 is "diff").
$oldHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$cols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $oldHeaders[$i] + "_FV2310"
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $oldHeaders[$i] + "_FV2404"
}

# 2) Freeze the header row (split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the data range into an Excel Table ("Table1").
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U91"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
